# CopperA-HW20 notebook rerun:
#  - two new simulation rows ("Holden", "Rizzie Spiral") inserted right after
#    the "Spiral5" row (i.e. they become the new rows 4 and 5), pushing all
#    the existing rows down by two.
#  - "Thomas Hex" renamed to "Matthies Hex".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shift the existing data block (rows 4-29, cols A-W) down by two
#     rows, to rows 6-31. Use Copy/PasteSpecial so formatting (and the
#     styles.xml table) is left exactly as-is rather than synthesizing new
#     style entries the way Rows.Insert() would.
$src = $ws.Range("A4:W29")
$dst = $ws.Range("A6:W31")
$src.Copy()
$dst.PasteSpecial(-4104)
$excel.CutCopyMode = 0

# --- Step 2: write the two new rows (4 and 5) with their labels and the
#     freshly recomputed simulation values.

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "Holden"

$row4Values = @(
    1.064081375591077,
    1.064081375591077,
    0.9935062546122492,
    1.017852990332688,
    1.030887920925754,
    0.9566676564085901,
    0.9267831034673433,
    0.9784473422804345,
    0.9925154148963079,
    0.9267831034673433,
    1.064081375591077,
    1.064081375591077,
    0.9925154148963079,
    0.9596492591818255,
    0.9930108347542785,
    0.994459964651576,
    0.9709349243253,
    0.994459964651576,
    0.9942215371417443,
    1.008193504831611,
    0.9950927573143056
)
$col = 3
foreach ($v in $row4Values) {
    $ws.Cells.Item(4, $col).Value = $v
    $col = $col + 1
}

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "Rizzie Spiral"

$row5Values = @(
    1.40056600203519,
    1.40056600203519,
    0.9355151632165263,
    1.026165493438321,
    1.091641389853606,
    0.9271796133506524,
    1.006839268817642,
    0.8513589182185937,
    0.8685884042816097,
    1.006839268817642,
    1.40056600203519,
    1.40056600203519,
    0.8685884042816097,
    0.9377138365496258,
    0.902051783749068,
    1.091997891711481,
    0.9369809454385926,
    1.091997891711481,
    1.052877209587742,
    1.122414968077231,
    1.013481781651518
)
$col = 3
foreach ($v in $row5Values) {
    $ws.Cells.Item(5, $col).Value = $v
    $col = $col + 1
}

# Copy the formatting (bold / bordered "index" style) from the row above
# onto the A/B cells of the two new rows, matching the rest of the column.
$ws.Range("A3:B3").Copy()
$ws.Range("A4:B5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 3: rename "Thomas Hex" -> "Matthies Hex". After the shift above,
#     that row (originally row 9) now lives at row 11.
$ws.Cells.Item(11, 2).Value = "Matthies Hex"

# --- Step 4: make sure the sheet's used range / dimension reflects the two
#     extra rows.
$ws.Range("A1:W31").Select()
